$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("ランサーズ")

$newTimestamp = "2025-12-18 06:38:52"

for ($r = 2; $r -le 10; $r++) {
    $ws.Cells.Item($r, 1).Value = $newTimestamp
}
